# Recomputed transition-probability matrix on Sheet1 after simulating more
# games (additional games added to the underlying simulation -> the
# per-state transition percentages shift slightly). Write the refreshed
# probabilities directly into the affected cells; columns/rows left at 0
# or otherwise untouched in the source diff are left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.213375796178344
$ws.Range("C2").Value = 0.5509554140127388
$ws.Range("J2").Value = 0.01910828025477707
$ws.Range("P2").Value = 0.143312101910828
$ws.Range("S2").Value = 0.0732484076433121
$ws.Range("B3").Value = 0.005524861878453038
$ws.Range("C3").Value = 0.03867403314917127
$ws.Range("J3").Value = 0.02209944751381215
$ws.Range("P3").Value = 0.8342541436464088
$ws.Range("S3").Value = 0.09944751381215469
$ws.Range("J4").Value = 0.04081632653061224
$ws.Range("P4").Value = 0.6938775510204082
$ws.Range("S4").Value = 0.2653061224489796
$ws.Range("B6").Value = 0.05445544554455446
$ws.Range("D6").Value = 0.0198019801980198
$ws.Range("F6").Value = 0.04455445544554455
$ws.Range("J6").Value = 0.2475247524752475
$ws.Range("O6").Value = 0.009900990099009901
$ws.Range("Q6").Value = 0.1683168316831683
$ws.Range("R6").Value = 0.07425742574257425
$ws.Range("S6").Value = 0.3811881188118812
$ws.Range("B7").Value = 0.1133333333333333
$ws.Range("D7").Value = 0.04
$ws.Range("F7").Value = 0.02666666666666667
$ws.Range("J7").Value = 0.08
$ws.Range("O7").Value = 0.02666666666666667
$ws.Range("Q7").Value = 0.18
$ws.Range("R7").Value = 0.1
$ws.Range("S7").Value = 0.4333333333333333
$ws.Range("B8").Value = 0.08917197452229299
$ws.Range("D8").Value = 0.02760084925690021
$ws.Range("F8").Value = 0.05307855626326964
$ws.Range("J8").Value = 0.09341825902335456
$ws.Range("O8").Value = 0.02123142250530785
$ws.Range("Q8").Value = 0.178343949044586
$ws.Range("R8").Value = 0.1125265392781316
$ws.Range("S8").Value = 0.4246284501061571
$ws.Range("B9").Value = 0.07983193277310924
$ws.Range("D9").Value = 0.008403361344537815
$ws.Range("E9").Value = 0.004201680672268907
$ws.Range("F9").Value = 0.07142857142857142
$ws.Range("J9").Value = 0.1176470588235294
$ws.Range("O9").Value = 0.01680672268907563
$ws.Range("Q9").Value = 0.1848739495798319
$ws.Range("R9").Value = 0.09243697478991597
$ws.Range("S9").Value = 0.4243697478991597
$ws.Range("B10").Value = 0.113933236574746
$ws.Range("D10").Value = 0.0181422351233672
$ws.Range("E10").Value = 0.001451378809869376
$ws.Range("F10").Value = 0.06168359941944847
$ws.Range("J10").Value = 0.1197387518142235
$ws.Range("O10").Value = 0.009433962264150943
$ws.Range("Q10").Value = 0.2155297532656023
$ws.Range("R10").Value = 0.08998548621190131
$ws.Range("S10").Value = 0.3701015965166908
$ws.Range("G11").Value = 0.1396396396396396
$ws.Range("J11").Value = 0.09009009009009009
$ws.Range("K11").Value = 0.1846846846846847
$ws.Range("L11").Value = 0.5765765765765766
$ws.Range("S11").Value = 0.009009009009009009
$ws.Range("G12").Value = 0.6985294117647058
$ws.Range("J12").Value = 0.1911764705882353
$ws.Range("K12").Value = 0.02941176470588235
$ws.Range("L12").Value = 0.04411764705882353
$ws.Range("S12").Value = 0.03676470588235294
$ws.Range("G13").Value = 0.717948717948718
$ws.Range("J13").Value = 0.2564102564102564
$ws.Range("S13").Value = 0.02564102564102564
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.02127659574468085
$ws.Range("H15").Value = 0.1702127659574468
$ws.Range("I15").Value = 0.1063829787234043
$ws.Range("J15").Value = 0.3563829787234042
$ws.Range("K15").Value = 0.04787234042553191
$ws.Range("M15").Value = 0.02659574468085106
$ws.Range("O15").Value = 0.03723404255319149
$ws.Range("S15").Value = 0.2340425531914894
$ws.Range("F16").Value = 0.008968609865470852
$ws.Range("H16").Value = 0.1569506726457399
$ws.Range("I16").Value = 0.08968609865470852
$ws.Range("J16").Value = 0.5112107623318386
$ws.Range("K16").Value = 0.04484304932735426
$ws.Range("M16").Value = 0.004484304932735426
$ws.Range("N16").Value = 0.004484304932735426
$ws.Range("O16").Value = 0.04484304932735426
$ws.Range("S16").Value = 0.1345291479820628
$ws.Range("F17").Value = 0.02092050209205021
$ws.Range("H17").Value = 0.1799163179916318
$ws.Range("I17").Value = 0.1317991631799163
$ws.Range("J17").Value = 0.4497907949790795
$ws.Range("K17").Value = 0.04811715481171548
$ws.Range("M17").Value = 0.01464435146443515
$ws.Range("O17").Value = 0.03765690376569038
$ws.Range("S17").Value = 0.1171548117154812
$ws.Range("F18").Value = 0.01327433628318584
$ws.Range("H18").Value = 0.1858407079646018
$ws.Range("I18").Value = 0.1150442477876106
$ws.Range("J18").Value = 0.4292035398230089
$ws.Range("K18").Value = 0.084070796460177
$ws.Range("M18").Value = 0.008849557522123894
$ws.Range("O18").Value = 0.05309734513274336
$ws.Range("S18").Value = 0.1106194690265487
$ws.Range("F19").Value = 0.01448170731707317
$ws.Range("H19").Value = 0.211890243902439
$ws.Range("I19").Value = 0.08231707317073171
$ws.Range("J19").Value = 0.4047256097560976
$ws.Range("K19").Value = 0.08689024390243902
$ws.Range("M19").Value = 0.01829268292682927
$ws.Range("N19").Value = 0.0007621951219512195
$ws.Range("O19").Value = 0.0625
$ws.Range("S19").Value = 0.118140243902439
